$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header column E (Assignee) - reuse the header style from D1 ---
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "Assignee"

# --- EM values for existing rows ---
# Row 2 reuses the "left aligned" style used by the rest of row 2 (A2,B2,D2)
$ws.Range("D2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("E2").Value = "EM"

# Rows 3-6 get the plain default-styled value (no formatting applied)
$ws.Range("E3").Value = "EM"
$ws.Range("E4").Value = "EM"
$ws.Range("E5").Value = "EM"
$ws.Range("E6").Value = "EM"

# --- New row 7: "selected download options" issue ---
$ws.Range("A4").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = "Sprint 3"

$ws.Range("B4").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("B7").Value = "selected download options"

$ws.Range("C7").Value = "The selected download options are not set when opening the STO widget"
$ws.Range("C7").WrapText = $true

$ws.Range("D4").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("D7").Value = "Open"

$ws.Range("E7").Value = "EM"

$ws.Rows.Item(7).RowHeight = 45

$ws.Application.CutCopyMode = $false

$ws.Range("E7").Select()
